$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-91 down to 79-92.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new "Camote" price record.
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44736
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = 100114002
$ws.Range("G78").Value = "Camote"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 20
$ws.Range("K78").Value = 20000
$ws.Range("L78").Value = 20000
$ws.Range("M78").Value = 20000
$ws.Range("N78").Value = "$/malla 20 kilos"
$ws.Range("O78").Value = "Perú"
$ws.Range("P78").Value = 1000
$ws.Range("Q78").Value = 20
$ws.Range("R78").Value = "Hortaliza"
